$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row20 = @(0.19, 4427.371668214, 35.672811882, 4391.698856334, 25.31023158866667, 1708.110787172, 22.116897024, 1685.993890148, 26.56685283933334)
$row21 = @(0.2, 4719.73852585, 35.270357802, 4684.468168048, 25.32196972466667, 1798.145568474, 21.926013394, 1776.21955508, 26.58711745933334)

for ($i = 0; $i -lt $row20.Length; $i++) {
    $ws.Cells.Item(20, $i + 1).Value = $row20[$i]
}

for ($i = 0; $i -lt $row21.Length; $i++) {
    $ws.Cells.Item(21, $i + 1).Value = $row21[$i]
}
